$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update numeric columns E:T for existing rows 2-4 (cluster labels A-D unchanged)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1985.706367333333
$ws.Range("H2").Value = 5957.119102000001
$ws.Range("I2").Value = 0.9998048086715072
$ws.Range("J2").Value = 0.9998048086715072
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 22.48784766666667
$ws.Range("N2").Value = 67.463543
$ws.Range("O2").Value = 0.4520839499795984
$ws.Range("P2").Value = 0.4520839499795983
$ws.Range("Q2").Value = 44654.26229932204
$ws.Range("R2").Value = 401888.3606938985
$ws.Range("S2").Value = 0.4519957071128116
$ws.Range("T2").Value = 0.4519957071128115

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1985.706367333333
$ws.Range("H3").Value = 5957.119102000001
$ws.Range("I3").Value = 0.9998048086715072
$ws.Range("J3").Value = 0.9998048086715072
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 23.26810333333333
$ws.Range("N3").Value = 69.80431
$ws.Range("O3").Value = 0.4677698025791556
$ws.Range("P3").Value = 0.4677698025791556
$ws.Range("Q3").Value = 46203.62094476996
$ws.Range("R3").Value = 415832.5885029297
$ws.Range("S3").Value = 0.4676784979699614
$ws.Range("T3").Value = 0.4676784979699614

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1985.706367333333
$ws.Range("H4").Value = 5957.119102000001
$ws.Range("I4").Value = 0.9998048086715072
$ws.Range("J4").Value = 0.9998048086715072
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.986685666666666
$ws.Range("N4").Value = 11.960057
$ws.Range("O4").Value = 0.08014624744124609
$ws.Range("P4").Value = 0.08014624744124607
$ws.Range("Q4").Value = 7916.387112856534
$ws.Range("R4").Value = 71247.48401570882
$ws.Range("S4").Value = 0.08013060358873432
$ws.Range("T4").Value = 0.0801306035887343

# Add new rows 5-7 for FAPs sending-cluster pairs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "S100a8"
$ws.Range("C5").Value = "Tlr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3876683333333333
$ws.Range("H5").Value = 1.163005
$ws.Range("I5").Value = 0.0001951913284927648
$ws.Range("J5").Value = 0.0001951913284927648
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 22.48784766666667
$ws.Range("N5").Value = 67.463543
$ws.Range("O5").Value = 0.4520839499795984
$ws.Range("P5").Value = 0.4520839499795983
$ws.Range("Q5").Value = 8.717826425190555
$ws.Range("R5").Value = 78.46043782671501
$ws.Range("S5").Value = 0.00008824286678677445
$ws.Range("T5").Value = 0.00008824286678677442

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "S100a8"
$ws.Range("C6").Value = "Tlr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3876683333333333
$ws.Range("H6").Value = 1.163005
$ws.Range("I6").Value = 0.0001951913284927648
$ws.Range("J6").Value = 0.0001951913284927648
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 23.26810333333333
$ws.Range("N6").Value = 69.80431
$ws.Range("O6").Value = 0.4677698025791556
$ws.Range("P6").Value = 0.4677698025791556
$ws.Range("Q6").Value = 9.02030683906111
$ws.Range("R6").Value = 81.18276155155
$ws.Range("S6").Value = 0.00009130460919422372
$ws.Range("T6").Value = 0.00009130460919422369

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "S100a8"
$ws.Range("C7").Value = "Tlr4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.3876683333333333
$ws.Range("H7").Value = 1.163005
$ws.Range("I7").Value = 0.0001951913284927648
$ws.Range("J7").Value = 0.0001951913284927648
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.986685666666666
$ws.Range("N7").Value = 11.960057
$ws.Range("O7").Value = 0.08014624744124609
$ws.Range("P7").Value = 0.08014624744124607
$ws.Range("Q7").Value = 1.545511787920556
$ws.Range("R7").Value = 13.909606091285
$ws.Range("S7").Value = 0.00001564385251176668
$ws.Range("T7").Value = 0.00001564385251176667

